$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Generator Data")
$ws1.Range("B2").Value = 12948.8666456
$ws1.Range("B3").Value = 5439.81887781656
$ws1.Range("B4").Value = 543.9818877816559
$ws1.Range("B5").Value = 89353.44855279999

$ws2 = $wb.Worksheets.Item("Yearly Fuel Costs")
$ws2.Range("B2").Value = 24439.53923883097
$ws2.Range("B3").Value = 24443.11044667862
$ws2.Range("B4").Value = 24443.1310154965
$ws2.Range("B5").Value = 24443.15038968634
